$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("arbolu")
$ws.Range("C2").Value = 0.5686274509803922
$ws.Range("D2").Value = 0.09998077477632913
$ws.Range("E2").Value = 0.7058823529411765
$ws.Range("F2").Value = 0.6602434077079107
$ws.Range("G2").Value = 0.7058823529411765
$ws.Range("H2").Value = 0.6384803921568627
$ws.Range("I2").Value = 0.8295847750865052

$ws = $wb.Worksheets.Item("bosqueu")
$ws.Range("A2").Value = 111
$ws.Range("C2").Value = 0.6666666666666666
$ws.Range("D2").Value = 0.04999038738816457
$ws.Range("E2").Value = 0.6764705882352942
$ws.Range("F2").Value = 0.7005347593582888
$ws.Range("G2").Value = 0.6764705882352942
$ws.Range("H2").Value = 0.666547106647537
$ws.Range("I2").Value = 0.8760092272202998

$ws = $wb.Worksheets.Item("knnu")
$ws.Range("C2").Value = 0.5588235294117647
$ws.Range("E2").Value = 0.5588235294117647
$ws.Range("F2").Value = 0.3122837370242215
$ws.Range("G2").Value = 0.5588235294117647
$ws.Range("H2").Value = 0.4006659267480577
$ws.Range("I2").Value = 0.8166089965397924

$ws = $wb.Worksheets.Item("arbolts")
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 73
$ws.Range("C2").Value = 0.5366620114553465
$ws.Range("D2").Value = 0.4870547248066527
$ws.Range("E2").Value = 0.6978930611538223
$ws.Range("F2").Value = 0.1796952182466954

$ws = $wb.Worksheets.Item("bosquets")
$ws.Range("C2").Value = 0.3900736969775532
$ws.Range("D2").Value = 0.3016540288647431
$ws.Range("E2").Value = 0.5492303968870833
$ws.Range("F2").Value = 0.2828324552765196

$ws = $wb.Worksheets.Item("knnts")
$ws.Range("C2").Value = 0.4210041135149655
$ws.Range("D2").Value = 0.3372007051283547
$ws.Range("E2").Value = 0.5806898527857661
$ws.Range("F2").Value = 0.1983219893132574

$ws = $wb.Worksheets.Item("arboltd")
$ws.Range("A2").Value = 16
$ws.Range("B2").Value = 141
$ws.Range("C2").Value = 0.6371500747583697
$ws.Range("D2").Value = 0.5827834163807639
$ws.Range("E2").Value = 0.7634025257888292
$ws.Range("F2").Value = 0.06027717482050421

$ws = $wb.Worksheets.Item("bosquetd")
$ws.Range("A2").Value = 119
$ws.Range("B2").Value = 141
$ws.Range("C2").Value = 0.4649502752480864
$ws.Range("D2").Value = 0.3627717795584061
$ws.Range("E2").Value = 0.6023053872898748
$ws.Range("F2").Value = 0.415040112673201

$ws = $wb.Worksheets.Item("knntd")
$ws.Range("B2").Value = 38
$ws.Range("C2").Value = 0.7350670775892014
$ws.Range("D2").Value = 0.7878212258727322
$ws.Range("E2").Value = 0.8875929392873358
$ws.Range("F2").Value = 0.2196529674477486

$ws = $wb.Worksheets.Item("arbolcc")
$ws.Range("A2").Value = 8
$ws.Range("C2").Value = 0.5954445053546944
$ws.Range("D2").Value = 0.8330168955623509
$ws.Range("E2").Value = 0.9126975926134302
$ws.Range("F2").Value = 0.2402945822555412

$ws = $wb.Worksheets.Item("bosquecc")
$ws.Range("A2").Value = 109
$ws.Range("C2").Value = 0.456863285031173
$ws.Range("D2").Value = 0.5143437335846587
$ws.Range("E2").Value = 0.7171776164832939
$ws.Range("F2").Value = 0.6299190801888896

$ws = $wb.Worksheets.Item("knncc")
$ws.Range("A2").Value = 7
$ws.Range("C2").Value = 0.6175794785234507
$ws.Range("D2").Value = 0.7743973631321944
$ws.Range("E2").Value = 0.8799985017783806
$ws.Range("F2").Value = 0.2937551742436224

$ws = $wb.Worksheets.Item("bosquepp")
$ws.Range("C2").Value = 0.5609322990703484
$ws.Range("D2").Value = 0.5091111489312204
$ws.Range("E2").Value = 0.7135202512411406
$ws.Range("F2").Value = 0.3669522732597754

$ws = $wb.Worksheets.Item("knnpp")
$ws.Range("C2").Value = 0.7417394081863822
$ws.Range("D2").Value = 1.077288484588244
$ws.Range("E2").Value = 1.037925086212027
$ws.Range("F2").Value = 0.1720426223806321
